$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, matching the source data
# (values like "517.91" or "1.00" or "58.151.68" must not be reinterpreted as numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.151.68'
$ws.Range("E2").Value = '  -4.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.607.84'
$ws.Range("E3").Value = '  -3.62%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '517.91'
$ws.Range("E5").Value = '  -1.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.40'
$ws.Range("E6").Value = '  -2.23%  '

$ws.Range("E7").Value = '  +0.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -1.60%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.71'
$ws.Range("E9").Value = '  -1.02%  '

$ws.Range("E10").Value = '  -2.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.338'
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.065.05'
$ws.Range("E13").Value = '  -3.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.127.19'
$ws.Range("E14").Value = '  -4.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.92'
$ws.Range("E15").Value = '  -1.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000136'
$ws.Range("E16").Value = '  -1.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.627.48'
$ws.Range("E17").Value = '  -12.50%  '

$ws.Range("E18").Value = '  -2.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '335.04'
$ws.Range("E19").Value = '  -2.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.38'
$ws.Range("E20").Value = '  -2.37%  '

$ws.Range("E21").Value = '  -3.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.12'
$ws.Range("E23").Value = '  +1.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.415'
$ws.Range("E24").Value = '  -1.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  -2.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.11'
$ws.Range("E27").Value = '  -2.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0789'
$ws.Range("E28").Value = '  -3.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.63'
$ws.Range("E29").Value = '  -3.06%  '

$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.58'
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.88'
$ws.Range("E32").Value = '  +0.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.74'
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.09'
$ws.Range("E34").Value = '  -4.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.18'
$ws.Range("E35").Value = '  -4.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.894'
$ws.Range("E36").Value = '  -4.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.845'
$ws.Range("E37").Value = '  -3.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.18'
$ws.Range("E38").Value = '  -2.41%  '

$ws.Range("E39").Value = '  -5.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.62'
$ws.Range("E40").Value = '  -1.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.44%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.598'
$ws.Range("E42").Value = '  -2.07%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0965'
$ws.Range("E43").Value = '  -2.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '268.43'
$ws.Range("E44").Value = '  -4.85%  '

$ws.Range("E45").Value = '  +0.78%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.10'
$ws.Range("E46").Value = '  -4.91%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0532'
$ws.Range("E47").Value = '  -1.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.041.60'
$ws.Range("E48").Value = '  -4.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0229'
$ws.Range("E49").Value = '  -1.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.64'
$ws.Range("E50").Value = '  -5.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.18'
$ws.Range("E51").Value = '  -4.80%  '
